$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename "managed store" (매점) API section labels to "branch" (지점) wording ---
$ws.Range("E18").Value = "지점 정보 생성"
$ws.Range("E19").Value = "모든 지점 정보"
$ws.Range("E20").Value = "한 지점 정보"
$ws.Range("E21").Value = "지점 정보 수정"
$ws.Range("E22").Value = "지점 정보 삭제"

# --- Rename the statistics function names / URI to the new, more explicit names ---
$ws.Range("F45").Value = "getAllSalesConsumerGroupByConsumerGender"
$ws.Range("F46").Value = "getAllSalesConsumerGroupByConsumerAge"
$ws.Range("C47").Value = "/statistics/sales_hour/all"
$ws.Range("F47").Value = "getAllSalesConsumerGroupBySalesHour"

# --- Add the new "per item" statistics API rows ---
$ws.Range("A48").Value = "SCM"
$ws.Range("B48").Value = "통계"
$ws.Range("C48").Value = "/statistics/consumer_gender/one_item/{find_item_num}"
$ws.Range("D48").Value = "GET"
$ws.Range("E48").Value = "특정 물품에 대한 성별 판매량"
$ws.Range("F48").Value = "getOneItemGroupByConsumerGender"
$ws.Range("G48").Value = "ok"

$ws.Range("A49").Value = "SCM"
$ws.Range("B49").Value = "통계"
$ws.Range("C49").Value = "/statistics/consumer_age/one_item/{find_item_num}"
$ws.Range("D49").Value = "GET"
$ws.Range("E49").Value = "특정 물품에 대한 나이대별 판매량"
$ws.Range("F49").Value = "getOneItemGroupByConsumerAge"
$ws.Range("G49").Value = "ok"

$ws.Range("A50").Value = "SCM"
$ws.Range("B50").Value = "통계"
$ws.Range("C50").Value = "/statistics/sales_hour/one_item/{find_item_num}"
$ws.Range("D50").Value = "GET"
$ws.Range("E50").Value = "특정 물품에 대한 시간대별 판매량"
$ws.Range("F50").Value = "getOneItemGroupBySalesHour"
$ws.Range("G50").Value = "ok"

# --- Restore view state (active cell moved down as new rows were authored) ---
$ws.Range("F51").Select()
